# Applies the cryptos list update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.099.47'
$ws.Range('E2').Value = '  -0.66%  '
$ws.Range('D3').Value = '3.271.19'
$ws.Range('E3').Value = '  +0.49%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''599.42'
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('D6').Value = '''138.10'
$ws.Range('E6').Value = '  -2.48%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.269.58'
$ws.Range('E8').Value = '  +0.65%  '
$ws.Range('D9').Value = '''0.511'
$ws.Range('E9').Value = '  -1.47%  '
$ws.Range('E10').Value = '  -0.28%  '
$ws.Range('D11').Value = '''5.46'
$ws.Range('E11').Value = '  +0.79%  '
$ws.Range('D12').Value = '''0.462'
$ws.Range('E12').Value = '  -1.06%  '
$ws.Range('E13').Value = '  -2.47%  '
$ws.Range('D14').Value = '''34.02'
$ws.Range('E14').Value = '  -1.25%  '
$ws.Range('D15').Value = '3.808.79'
$ws.Range('E15').Value = '  +0.47%  '
$ws.Range('D16').Value = '''0.122'
$ws.Range('E16').Value = '  +1.30%  '
$ws.Range('D17').Value = '3.269.63'
$ws.Range('E17').Value = '  +0.53%  '
$ws.Range('D18').Value = '63.144.14'
$ws.Range('E18').Value = '  -0.62%  '
$ws.Range('D19').Value = '''6.75'
$ws.Range('E19').Value = '  -0.63%  '
$ws.Range('D20').Value = '''472.15'
$ws.Range('E20').Value = '  -1.34%  '
$ws.Range('D21').Value = '''13.82'
$ws.Range('E21').Value = '  -2.99%  '
$ws.Range('E22').Value = '  -0.88%  '
$ws.Range('D23').Value = '''7.86'
$ws.Range('E23').Value = '  -1.64%  '
$ws.Range('D24').Value = '''13.67'
$ws.Range('E24').Value = '  +2.83%  '
$ws.Range('D25').Value = '''84.21'
$ws.Range('E25').Value = '  +0.42%  '
$ws.Range('E27').Value = '  -0.70%  '
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('E29').Value = '  -1.93%  '
$ws.Range('E30').Value = '  -1.74%  '
$ws.Range('E31').Value = '  -1.40%  '
$ws.Range('D32').Value = '''28.14'
$ws.Range('E32').Value = '  +1.51%  '
$ws.Range('E33').Value = '  -2.83%  '
$ws.Range('D34').Value = '''2.47'
$ws.Range('E34').Value = '  -3.44%  '
$ws.Range('E35').Value = '  -1.10%  '
$ws.Range('E36').Value = '  -0.36%  '
$ws.Range('D37').Value = '''51.74'
$ws.Range('E37').Value = '  -1.92%  '
$ws.Range('D38').Value = '0.0₃0720'
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('D39').Value = '''0.0394'
$ws.Range('E39').Value = '  -0.20%  '
$ws.Range('D40').Value = '3.081.84'
$ws.Range('E40').Value = '  +2.70%  '
$ws.Range('D41').Value = '''422.73'
$ws.Range('E41').Value = '  -0.12%  '
$ws.Range('E42').Value = '  +6.05%  '
$ws.Range('E43').Value = '  -2.12%  '
$ws.Range('D44').Value = '''2.65'
$ws.Range('E44').Value = '  -4.44%  '
$ws.Range('D45').Value = '''0.258'
$ws.Range('E45').Value = '  -3.02%  '
$ws.Range('E46').Value = '  -0.98%  '
$ws.Range('B48').Value = 'Arweave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D48').Value = '''35.88'
$ws.Range('E48').Value = '  +6.50%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').Value = '''127.15'
$ws.Range('E49').Value = '  +3.51%  '
$ws.Range('D50').Value = '''25.85'
$ws.Range('E50').Value = '  -0.44%  '
$ws.Range('E51').Value = '  -1.66%  '

# The apostrophe above leaves a quote-prefix style on these cells;
# clear it so formatting matches the rest of the untouched column.
$ws.Range('D5').ClearFormats()
$ws.Range('D6').ClearFormats()
$ws.Range('D9').ClearFormats()
$ws.Range('D11').ClearFormats()
$ws.Range('D12').ClearFormats()
$ws.Range('D14').ClearFormats()
$ws.Range('D16').ClearFormats()
$ws.Range('D19').ClearFormats()
$ws.Range('D20').ClearFormats()
$ws.Range('D21').ClearFormats()
$ws.Range('D23').ClearFormats()
$ws.Range('D24').ClearFormats()
$ws.Range('D25').ClearFormats()
$ws.Range('D32').ClearFormats()
$ws.Range('D34').ClearFormats()
$ws.Range('D37').ClearFormats()
$ws.Range('D39').ClearFormats()
$ws.Range('D41').ClearFormats()
$ws.Range('D44').ClearFormats()
$ws.Range('D45').ClearFormats()
$ws.Range('D48').ClearFormats()
$ws.Range('D49').ClearFormats()
$ws.Range('D50').ClearFormats()
